$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: updated odds (existing row, selective cell updates) ---
$ws.Range("G2").Value = 1.75
$ws.Range("H2").Value = 4.2
$ws.Range("I2").Value = 4.1
$ws.Range("J2").Value = 2.25
$ws.Range("M2").Value = 1.02
$ws.Range("N2").Value = 19
$ws.Range("O2").Value = 1.14
$ws.Range("P2").Value = 5.5
$ws.Range("Q2").Value = 1.44
$ws.Range("R2").Value = 2.75
$ws.Range("W2").Value = 12
$ws.Range("X2").Value = 11
$ws.Range("Y2").Value = 8.5
$ws.Range("AC2").Value = 21
$ws.Range("AD2").Value = 8.5
$ws.Range("AF2").Value = 34
$ws.Range("AG2").Value = 101
$ws.Range("AH2").Value = 19
$ws.Range("AP2").Value = 15
$ws.Range("AR2").Value = 41
$ws.Range("AS2").Value = 81
$ws.Range("BB2").Value = 67

# --- Row 3: new row (Antalyaspor vs Bodrumspor) ---
$ws.Range("A3").Value = "WnPgMvn1"
$ws.Range("B3").Value = "'09/11/2024"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = "'07:30"
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").Value = "TURKEY - SUPER LIG"
$ws.Range("E3").Value = "Antalyaspor"
$ws.Range("F3").Value = "Bodrumspor"
$ws.Range("G3").Value = 2.3
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = 3.3
$ws.Range("J3").Value = 3.1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 3.75
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 8
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 2.75
$ws.Range("Q3").Value = 2.25
$ws.Range("R3").Value = 1.62
$ws.Range("S3").Value = 1.5
$ws.Range("T3").Value = 2.5
$ws.Range("U3").Value = 1.91
$ws.Range("V3").Value = 1.8
$ws.Range("W3").Value = 7
$ws.Range("X3").Value = 10
$ws.Range("Y3").Value = 10
$ws.Range("Z3").Value = 21
$ws.Range("AA3").Value = 21
$ws.Range("AB3").Value = 34
$ws.Range("AC3").Value = 7.5
$ws.Range("AD3").Value = 6
$ws.Range("AE3").Value = 15
$ws.Range("AF3").Value = 51
$ws.Range("AG3").Value = 351
$ws.Range("AH3").Value = 9
$ws.Range("AI3").Value = 15
$ws.Range("AJ3").Value = 12
$ws.Range("AK3").Value = 34
$ws.Range("AL3").Value = 29
$ws.Range("AM3").Value = 41
$ws.Range("AN3").Value = 4.33
$ws.Range("AO3").Value = 13
$ws.Range("AP3").Value = 26
$ws.Range("AQ3").Value = 41
$ws.Range("AR3").Value = 67
$ws.Range("AS3").Value = 201
$ws.Range("AT3").Value = 2.5
$ws.Range("AU3").Value = 8.5
$ws.Range("AV3").Value = 67
$ws.Range("AW3").Value = 276
$ws.Range("AX3").Value = 5
$ws.Range("AY3").Value = 19
$ws.Range("AZ3").Value = 29
$ws.Range("BA3").Value = 67
$ws.Range("BB3").Value = 101
$ws.Range("BC3").Value = 251
$ws.Range("BD3").Value = 276

# --- Row 4: new row (Istanbulspor AS vs Adanaspor AS) ---
$ws.Range("A4").Value = "86gnZDvo"
$ws.Range("B4").Value = "'09/11/2024"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = "'07:30"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = "TURKEY - 1. LIG"
$ws.Range("E4").Value = "Istanbulspor AS"
$ws.Range("F4").Value = "Adanaspor AS"
$ws.Range("G4").Value = 1.36
$ws.Range("H4").Value = 4.75
$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 1.83
$ws.Range("K4").Value = 2.5
$ws.Range("L4").Value = 7
$ws.Range("M4").Value = 1.03
$ws.Range("N4").Value = 15
$ws.Range("O4").Value = 1.18
$ws.Range("P4").Value = 4.5
$ws.Range("Q4").Value = 1.6
$ws.Range("R4").Value = 2.3
$ws.Range("S4").Value = 1.29
$ws.Range("T4").Value = 3.5
$ws.Range("U4").Value = 1.83
$ws.Range("V4").Value = 1.83
$ws.Range("W4").Value = 8
$ws.Range("X4").Value = 7
$ws.Range("Y4").Value = 8.5
$ws.Range("Z4").Value = 9
$ws.Range("AA4").Value = 11
$ws.Range("AB4").Value = 23
$ws.Range("AC4").Value = 15
$ws.Range("AD4").Value = 9.5
$ws.Range("AE4").Value = 19
$ws.Range("AF4").Value = 51
$ws.Range("AG4").Value = 251
$ws.Range("AH4").Value = 21
$ws.Range("AI4").Value = 41
$ws.Range("AJ4").Value = 21
$ws.Range("AK4").Value = 81
$ws.Range("AL4").Value = 51
$ws.Range("AM4").Value = 51
$ws.Range("AN4").Value = 3.4
$ws.Range("AO4").Value = 6.5
$ws.Range("AP4").Value = 17
$ws.Range("AQ4").Value = 17
$ws.Range("AR4").Value = 41
$ws.Range("AS4").Value = 101
$ws.Range("AT4").Value = 3.5
$ws.Range("AU4").Value = 9
$ws.Range("AV4").Value = 51
$ws.Range("AW4").Value = 126
$ws.Range("AX4").Value = 8.5
$ws.Range("AY4").Value = 34
$ws.Range("AZ4").Value = 41
$ws.Range("BA4").Value = 151
$ws.Range("BB4").Value = 151
$ws.Range("BC4").Value = 251
$ws.Range("BD4").Value = 126
